$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047658907431494
$ws.Range("D2").Value = 1.048344837462129
$ws.Range("E2").Value = 1.058376802429152
$ws.Range("F2").Value = 1.066248942139921
$ws.Range("I2").Value = 1.041508688651584
$ws.Range("J2").Value = 1.052706800491214
$ws.Range("K2").Value = 1.051104942810012
$ws.Range("L2").Value = 1.06110920413722
$ws.Range("M2").Value = 1.068959999701391
$ws.Range("N2").Value = 1.054201764895069

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049435552865936
$ws.Range("D3").Value = 1.049720972898623
$ws.Range("E3").Value = 1.05993623212775
$ws.Range("F3").Value = 1.067909521328594
$ws.Range("I3").Value = 1.042027919677412
$ws.Range("J3").Value = 1.054127946448451
$ws.Range("K3").Value = 1.052291645408929
$ws.Range("L3").Value = 1.062480765871848
$ws.Range("M3").Value = 1.07043402732839
$ws.Range("N3").Value = 1.055624929042573

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050581628555122
$ws.Range("D4").Value = 1.050608216332054
$ws.Range("E4").Value = 1.06094101470865
$ws.Range("F4").Value = 1.068980197803043
$ws.Range("I4").Value = 1.042360759715498
$ws.Range("J4").Value = 1.055043728151591
$ws.Range("K4").Value = 1.053055791144956
$ws.Range("L4").Value = 1.063363500374314
$ws.Range("M4").Value = 1.071383523681383
$ws.Range("N4").Value = 1.056542011260774

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051062611107193
$ws.Range("D5").Value = 1.050980458448197
$ws.Range("E5").Value = 1.06136242008602
$ws.Range("F5").Value = 1.069429411277858
$ws.Range("I5").Value = 1.042499940886447
$ws.Range("J5").Value = 1.055427829037891
$ws.Range("K5").Value = 1.053376157275092
$ws.Range("L5").Value = 1.06373347959474
$ws.Range("M5").Value = 1.071781679703469
$ws.Range("N5").Value = 1.056926657614389

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051143322180778
$ws.Range("D6").Value = 1.051042915638366
$ws.Range("E6").Value = 1.061433117466951
$ws.Range("F6").Value = 1.069504783998765
$ws.Range("I6").Value = 1.042523266523639
$ws.Range("J6").Value = 1.055492269282419
$ws.Range("K6").Value = 1.053429896865351
$ws.Range("L6").Value = 1.063795535397926
$ws.Range("M6").Value = 1.071848472975872
$ws.Range("N6").Value = 1.056991189371455

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050588058695107
$ws.Range("D7").Value = 1.050613193198361
$ws.Range("E7").Value = 1.060946649473736
$ws.Range("F7").Value = 1.068986203723366
$ws.Range("I7").Value = 1.042362622379033
$ws.Range("J7").Value = 1.055048864018155
$ws.Range("K7").Value = 1.053060075336805
$ws.Range("L7").Value = 1.063368448442302
$ws.Range("M7").Value = 1.071388847816331
$ws.Range("N7").Value = 1.056547154420857

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048260075475531
$ws.Range("D8").Value = 1.048810581042503
$ws.Range("E8").Value = 1.058904714378963
$ws.Range("F8").Value = 1.066810946315154
$ws.Range("I8").Value = 1.041684818566423
$ws.Range("J8").Value = 1.053187878064723
$ws.Range("K8").Value = 1.051506774302856
$ws.Range("L8").Value = 1.061573725005767
$ws.Range("M8").Value = 1.06945905430259
$ws.Range("N8").Value = 1.054683525653913

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044129999461962
$ws.Range("D9").Value = 1.045608988150957
$ws.Range("E9").Value = 1.055273042592385
$ws.Range("F9").Value = 1.062947732540468
$ws.Range("I9").Value = 1.040466131551278
$ws.Range("J9").Value = 1.049878854236686
$ws.Range("K9").Value = 1.048740535529212
$ws.Range("L9").Value = 1.058374004952552
$ws.Range("M9").Value = 1.066024838563213
$ws.Range("N9").Value = 1.05136980263249

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041356728396774
$ws.Range("D10").Value = 1.0434568472364
$ws.Range("E10").Value = 1.052828269884502
$ws.Range("F10").Value = 1.060350866815845
$ws.Range("I10").Value = 1.03963692928519
$ws.Range("J10").Value = 1.047651913715035
$ws.Range("K10").Value = 1.046876013239203
$ws.Range("L10").Value = 1.056214801191885
$ws.Range("M10").Value = 1.063711660649192
$ws.Range("N10").Value = 1.049139699599865

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040150897250525
$ws.Range("D11").Value = 1.04252055108532
$ws.Range("E11").Value = 1.051763794064084
$ws.Range("F11").Value = 1.059221068657623
$ws.Range("I11").Value = 1.039273814112397
$ws.Range("J11").Value = 1.046682450648252
$ws.Range("K11").Value = 1.046063649006397
$ws.Range("L11").Value = 1.055273425780503
$ws.Range("M11").Value = 1.062704169919472
$ws.Range("N11").Value = 1.048168859784328

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039702226413301
$ws.Range("D12").Value = 1.042172091019213
$ws.Range("E12").Value = 1.051367496249382
$ws.Range("F12").Value = 1.058800587377072
$ws.Range("I12").Value = 1.039138318383942
$ws.Range("J12").Value = 1.046321551943612
$ws.Range("K12").Value = 1.045761132224129
$ws.Range("L12").Value = 1.054922771163911
$ws.Range("M12").Value = 1.062329040587693
$ws.Range("N12").Value = 1.047807448562131

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039798503102269
$ws.Range("D13").Value = 1.042246867829006
$ws.Range("E13").Value = 1.051452544740146
$ws.Range("F13").Value = 1.058890819679438
$ws.Range("I13").Value = 1.039167410797114
$ws.Range("J13").Value = 1.046399002300436
$ws.Range("K13").Value = 1.045826058126138
$ws.Range("L13").Value = 1.054998032713459
$ws.Range("M13").Value = 1.062409548272581
$ws.Range("N13").Value = 1.047885008907346

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040113825846773
$ws.Range("D14").Value = 1.042491761192784
$ws.Range("E14").Value = 1.05173105453316
$ws.Range("F14").Value = 1.059186328479892
$ws.Range("I14").Value = 1.039262626657974
$ws.Range("J14").Value = 1.046652635012882
$ws.Range("K14").Value = 1.046038658641814
$ws.Range("L14").Value = 1.055244460791279
$ws.Range("M14").Value = 1.06267318014104
$ws.Range("N14").Value = 1.048139001807336

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040308003698487
$ws.Range("D15").Value = 1.042642557748465
$ws.Range("E15").Value = 1.051902533151793
$ws.Range("F15").Value = 1.059368291363468
$ws.Range("I15").Value = 1.039321210087081
$ws.Range("J15").Value = 1.046808800343254
$ws.Range("K15").Value = 1.046169546537238
$ws.Range("L15").Value = 1.055396162009028
$ws.Range("M15").Value = 1.062835492236332
$ws.Range("N15").Value = 1.048295388910389

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041436648482842
$ws.Range("D16").Value = 1.043518891944189
$ws.Range("E16").Value = 1.052898790050783
$ws.Range("F16").Value = 1.060425733463393
$ws.Range("I16").Value = 1.03966094172866
$ws.Range("J16").Value = 1.047716143031323
$ws.Range("K16").Value = 1.04692982024461
$ws.Range("L16").Value = 1.056277140050943
$ws.Range("M16").Value = 1.063778399076404
$ws.Range("N16").Value = 1.049204020129149

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04214326686817
$ws.Range("D17").Value = 1.044067402543828
$ws.Range("E17").Value = 1.053522128031078
$ws.Range("F17").Value = 1.061087594376778
$ws.Range("I17").Value = 1.039872952559011
$ws.Range("J17").Value = 1.048283894651643
$ws.Range("K17").Value = 1.047405366809213
$ws.Range("L17").Value = 1.056828019540893
$ws.Range("M17").Value = 1.06436827375341
$ws.Range("N17").Value = 1.049772578021892

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042554945815819
$ws.Range("D18").Value = 1.04438691558981
$ws.Range("E18").Value = 1.053885145340927
$ws.Range("F18").Value = 1.06147313309053
$ws.Range("I18").Value = 1.039996223200174
$ws.Range("J18").Value = 1.048614555714972
$ws.Range("K18").Value = 1.047682262360335
$ws.Range("L18").Value = 1.057148719611741
$ws.Range("M18").Value = 1.064711772762542
$ws.Range("N18").Value = 1.050103708661857

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042695237019783
$ws.Range("D19").Value = 1.044495789866202
$ws.Range("E19").Value = 1.054008829650195
$ws.Range("F19").Value = 1.061604505379885
$ws.Range("I19").Value = 1.040038189101624
$ws.Range("J19").Value = 1.048727218527179
$ws.Range("K19").Value = 1.04777659516978
$ws.Range("L19").Value = 1.057257965735341
$ws.Range("M19").Value = 1.064828801848719
$ws.Range("N19").Value = 1.050216531468182

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042067503170931
$ws.Range("D20").Value = 1.044008596510934
$ws.Range("E20").Value = 1.053455308381687
$ws.Range("F20").Value = 1.061016636242547
$ws.Range("I20").Value = 1.039850246374173
$ws.Range("J20").Value = 1.048223032005982
$ws.Range("K20").Value = 1.047354395202842
$ws.Range("L20").Value = 1.056768979518791
$ws.Range("M20").Value = 1.064305044362864
$ws.Range("N20").Value = 1.049711628944293

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040020992634353
$ws.Range("D21").Value = 1.042419665047899
$ws.Range("E21").Value = 1.051649065474261
$ws.Range("F21").Value = 1.059099331425935
$ws.Range("I21").Value = 1.039234605099709
$ws.Range("J21").Value = 1.046577968648389
$ws.Range("K21").Value = 1.045976074427926
$ws.Range("L21").Value = 1.055171921244896
$ws.Range("M21").Value = 1.062595572204508
$ws.Range("N21").Value = 1.04806422940804

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038729795543063
$ws.Range("D22").Value = 1.041416709962735
$ws.Range("E22").Value = 1.050508167062721
$ws.Range("F22").Value = 1.0578890676494
$ws.Range("I22").Value = 1.038843944081623
$ws.Range("J22").Value = 1.045539031858845
$ws.Range("K22").Value = 1.045105015221287
$ws.Range("L22").Value = 1.054162073038483
$ws.Range("M22").Value = 1.061515527708601
$ws.Range("N22").Value = 1.047023817209088

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039414715222571
$ws.Range("D23").Value = 1.04194877349079
$ws.Range("E23").Value = 1.051113483087569
$ws.Range("F23").Value = 1.05853111168852
$ws.Range("I23").Value = 1.039051383146745
$ws.Range("J23").Value = 1.046090236004156
$ws.Range("K23").Value = 1.045567207551279
$ws.Range("L23").Value = 1.054697961222927
$ws.Range("M23").Value = 1.062088582564361
$ws.Range("N23").Value = 1.047575804127509

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.04210173900592
$ws.Range("D24").Value = 1.044035169735466
$ws.Range("E24").Value = 1.053485503050853
$ws.Range("F24").Value = 1.061048700755366
$ws.Range("I24").Value = 1.039860507531859
$ws.Range("J24").Value = 1.048250534756767
$ws.Range("K24").Value = 1.047377428568953
$ws.Range("L24").Value = 1.056795659074353
$ws.Range("M24").Value = 1.064333616747778
$ws.Range("N24").Value = 1.049739170752138

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045201143340412
$ws.Range("D25").Value = 1.046439742961976
$ws.Range("E25").Value = 1.056216003765576
$ws.Range("F25").Value = 1.063950152289787
$ws.Range("I25").Value = 1.040784113016582
$ws.Range("J25").Value = 1.050737934915394
$ws.Range("K25").Value = 1.049459205740369
$ws.Range("L25").Value = 1.059205723990151
$ws.Range("M25").Value = 1.066916764504808
$ws.Range("N25").Value = 1.052230103304299
